# Refresh cached Universalis market-price figures (columns H:N) for the
# Leve-profit rows flagged by the scheduled pricing runner, one worksheet
# (Job abbreviation) at a time.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 929
$ws.Range("I98").Value = 540.6
$ws.Range("J98").Value = 1900
$ws.Range("K98").Value = 540.6
$ws.Range("L98").Value = 1900
$ws.Range("M98").Value = 957.4
$ws.Range("N98").Value = -4896

# Row 122
$ws.Range("H122").Value = 929
$ws.Range("I122").Value = 540.6
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 1621.8
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = 828.1999999999998
$ws.Range("N122").Value = -10600

# Row 138
$ws.Range("H138").Value = 1532.7606
$ws.Range("I138").Value = 1932.5834
$ws.Range("J138").Value = 1328.5957
$ws.Range("K138").Value = 5797.7502
$ws.Range("L138").Value = 3985.7871
$ws.Range("M138").Value = -657.7502000000004
$ws.Range("N138").Value = -14265.7871

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2924.5173
$ws.Range("I61").Value = 2161.7222
$ws.Range("J61").Value = 4172.727
$ws.Range("K61").Value = 2161.7222
$ws.Range("L61").Value = 4172.727
$ws.Range("M61").Value = -1949.7222
$ws.Range("N61").Value = -4596.727

# Row 132
$ws.Range("H132").Value = 6856.364
$ws.Range("I132").Value = 14004
$ws.Range("J132").Value = 4176
$ws.Range("K132").Value = 42012
$ws.Range("L132").Value = 12528
$ws.Range("M132").Value = -39482
$ws.Range("N132").Value = -17588

# Row 136
$ws.Range("H136").Value = 2924.5173
$ws.Range("I136").Value = 2161.7222
$ws.Range("J136").Value = 4172.727
$ws.Range("K136").Value = 6485.1666
$ws.Range("L136").Value = 12518.181
$ws.Range("M136").Value = -3935.1666
$ws.Range("N136").Value = -17618.181

$ws = $wb.Worksheets.Item("BSM")
# Row 29
$ws.Range("H29").Value = 2508
$ws.Range("I29").Value = 2508
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 2508
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -2219

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1518.2142
$ws.Range("I58").Value = 1042.8572
$ws.Range("J58").Value = 1993.5714
$ws.Range("K58").Value = 1042.8572
$ws.Range("L58").Value = 1993.5714
$ws.Range("M58").Value = -839.8571999999999
$ws.Range("N58").Value = -2399.5714

# Row 92
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# Row 116
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

# Row 119
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

# Row 132
$ws.Range("H132").Value = 5748811
$ws.Range("I132").Value = 1093.1428
$ws.Range("J132").Value = 20836572
$ws.Range("K132").Value = 3279.4284
$ws.Range("L132").Value = 62509716
$ws.Range("M132").Value = -749.4284000000002
$ws.Range("N132").Value = -62514776

# Row 136
$ws.Range("H136").Value = 1518.2142
$ws.Range("I136").Value = 1042.8572
$ws.Range("J136").Value = 1993.5714
$ws.Range("K136").Value = 3128.5716
$ws.Range("L136").Value = 5980.7142
$ws.Range("M136").Value = -578.5715999999998
$ws.Range("N136").Value = -11080.7142

# Row 140
$ws.Range("H140").Value = 68800
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 68800
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 68800
$ws.Range("N140").Value = -79160

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 10002829
$ws.Range("I4").Value = 10003250
$ws.Range("J4").Value = 10002661
$ws.Range("K4").Value = 30009750
$ws.Range("L4").Value = 30007983
$ws.Range("M4").Value = -30009638
$ws.Range("N4").Value = -30008207

# Row 39
$ws.Range("H39").Value = 1390.4615
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 1786.2222
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 5358.6666
$ws.Range("M39").Value = -1206
$ws.Range("N39").Value = -5946.6666

# Row 55
$ws.Range("H55").Value = 1700
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 2300
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 6900
$ws.Range("M55").Value = -1323
$ws.Range("N55").Value = -7254

# Row 110
$ws.Range("H110").Value = 12742.7
$ws.Range("I110").Value = 5009
$ws.Range("J110").Value = 16057.143
$ws.Range("K110").Value = 15027
$ws.Range("L110").Value = 48171.429
$ws.Range("M110").Value = -10937
$ws.Range("N110").Value = -56351.429

# Row 134
$ws.Range("H134").Value = 5574.0835
$ws.Range("I134").Value = 4814.8335
$ws.Range("J134").Value = 6333.3335
$ws.Range("K134").Value = 14444.5005
$ws.Range("L134").Value = 19000.0005
$ws.Range("M134").Value = -9374.500499999998
$ws.Range("N134").Value = -29140.0005

# Row 139
$ws.Range("H139").Value = 3522.25
$ws.Range("I139").Value = 2908.4285
$ws.Range("J139").Value = 3999.6667
$ws.Range("K139").Value = 8725.2855
$ws.Range("L139").Value = 11999.0001
$ws.Range("M139").Value = -3585.2855
$ws.Range("N139").Value = -22279.0001

# Row 140
$ws.Range("H140").Value = 1119.7142
$ws.Range("I140").Value = 1025.7
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 3077.1
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = 2102.9
$ws.Range("N140").Value = -19360

# Row 141
$ws.Range("H141").Value = 4372
$ws.Range("I141").Value = 1477.375
$ws.Range("J141").Value = 7266.625
$ws.Range("K141").Value = 4432.125
$ws.Range("L141").Value = 21799.875
$ws.Range("M141").Value = 747.875
$ws.Range("N141").Value = -32159.875

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2806.9375
$ws.Range("I126").Value = 2291.8
$ws.Range("J126").Value = 3665.5
$ws.Range("K126").Value = 6875.400000000001
$ws.Range("L126").Value = 10996.5
$ws.Range("M126").Value = -4405.400000000001
$ws.Range("N126").Value = -15936.5

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 4826.9546
$ws.Range("I68").Value = 3758.6
$ws.Range("J68").Value = 5141.1763
$ws.Range("K68").Value = 3758.6
$ws.Range("L68").Value = 5141.1763
$ws.Range("M68").Value = -3009.6
$ws.Range("N68").Value = -6639.1763

# Row 71
$ws.Range("H71").Value = 4826.9546
$ws.Range("I71").Value = 3758.6
$ws.Range("J71").Value = 5141.1763
$ws.Range("K71").Value = 18793
$ws.Range("L71").Value = 25705.8815
$ws.Range("M71").Value = -15049
$ws.Range("N71").Value = -33193.8815

# Row 136
$ws.Range("H136").Value = 2414.65
$ws.Range("I136").Value = 2661
$ws.Range("J136").Value = 1957.1428
$ws.Range("K136").Value = 7983
$ws.Range("L136").Value = 5871.428400000001
$ws.Range("M136").Value = -5433
$ws.Range("N136").Value = -10971.4284

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 3310.6072
$ws.Range("I136").Value = 3047.9375
$ws.Range("J136").Value = 3660.8333
$ws.Range("K136").Value = 9143.8125
$ws.Range("L136").Value = 10982.4999
$ws.Range("M136").Value = -6593.8125
$ws.Range("N136").Value = -16082.4999
